$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "10815a"
$ws.Range("A3").Value = "11069b"
$ws.Range("A4").Value = "10815a"
$ws.Range("A5").Value = "11070a"
$ws.Range("A6").Value = "10815a"
$ws.Range("A7").Value = "11071a"
$ws.Range("A8").Value = "10815a"
$ws.Range("A9").Value = "11072a"
$ws.Range("A10").Value = "10815a"
$ws.Range("A11").Value = "11074b"
$ws.Range("A12").Value = "10815a"
